$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text looks like a plain number (e.g. "248.11") need the
# cell format forced to Text first, otherwise Excel auto-converts the entry
# into a numeric value instead of keeping it as the original text string.

$ws.Range('D2').Value = '34.751.70'
$ws.Range('E2').Value = '  -2.23%  '

$ws.Range('D3').Value = '1.875.20'
$ws.Range('E3').Value = '  -1.99%  '

$ws.Range('E4').Value = '  -0.94%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '248.11'
$ws.Range('E5').Value = '  +0.29%  '

$ws.Range('E6').Value = '  -2.16%  '

$ws.Range('E7').Value = '  -0.95%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '41.38'
$ws.Range('E8').Value = '  +1.39%  '

$ws.Range('E9').Value = '  -2.67%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '50.65'
$ws.Range('E10').Value = '  -4.20%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0740'
$ws.Range('E11').Value = '  +0.34%  '

$ws.Range('E12').Value = '  -2.36%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '12.87'
$ws.Range('E13').Value = '  +1.14%  '

$ws.Range('D14').Value = '2.148.40'
$ws.Range('E14').Value = '  -1.94%  '

$ws.Range('E15').Value = '  -0.50%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '4.90'
$ws.Range('E16').Value = '  -0.48%  '

$ws.Range('D17').Value = '1.866.02'
$ws.Range('E17').Value = '  -2.58%  '

$ws.Range('D18').Value = '34.735.84'
$ws.Range('E18').Value = '  -2.25%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '72.88'
$ws.Range('E19').Value = '  -0.59%  '

$ws.Range('E20').Value = '  -0.43%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '247.59'
$ws.Range('E21').Value = '  +2.01%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '12.76'
$ws.Range('E22').Value = '  -3.29%  '

$ws.Range('E23').Value = '  -3.35%  '

$ws.Range('E24').Value = '  -0.99%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.41'
$ws.Range('E25').Value = '  +3.74%  '

$ws.Range('E26').Value = '  -1.65%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '165.32'

$ws.Range('E28').Value = '  -2.98%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '18.25'
$ws.Range('E29').Value = '  -3.31%  '

$ws.Range('E30').Value = '  -3.67%  '

$ws.Range('D31').Value = '4.128.36'
$ws.Range('E31').Value = '  -0.34%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.71'
$ws.Range('E32').Value = '  +14.80%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.25'
$ws.Range('E33').Value = '  -0.30%  '

$ws.Range('E34').Value = '  +0.06%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.17'
$ws.Range('E35').Value = '  -1.48%  '

$ws.Range('E36').Value = '  -1.00%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.81'
$ws.Range('E37').Value = '  -5.95%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.836'
$ws.Range('E38').Value = '  -8.96%  '

$ws.Range('E39').Value = '  -3.54%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('B40').Value = 'InjectiveProtocol'
$ws.Range('C40').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D40').Value = '17.25'
$ws.Range('E40').Value = '  -2.75%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('B41').Value = 'Aave'
$ws.Range('C41').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D41').Value = '98.30'
$ws.Range('E41').Value = '  -0.72%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0657'
$ws.Range('E42').Value = '  +1.62%  '

$ws.Range('E43').Value = '  -0.06%  '

$ws.Range('E44').Value = '  -5.49%  '

$ws.Range('D45').Value = '1.293.77'
$ws.Range('E45').Value = '  -4.38%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.36'
$ws.Range('E46').Value = '  -4.52%  '

$ws.Range('E47').Value = '  -0.88%  '

$ws.Range('E48').Value = '  -2.24%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0766'
$ws.Range('E49').Value = '  +6.25%  '

$ws.Range('E50').Value = '  -1.31%  '

$ws.Range('E51').Value = '  -1.28%  '
